# Insert 4 new data rows at row 448 (pushing the existing rows 448:552 down
# to 452:556) and populate the newly inserted rows with the new daily
# price-report entries for 2022-10-?? (serial date 44855).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 448.
$ws.Rows.Item(448).Resize(4).Insert()

# Columns that are constant across every record in this sheet.
$constA = 4
$constB = "Feria Lagunitas de Puerto Montt"
$constC = "Los Lagos"
$constE = 10
$constF = 100112006
$constG = "Repollo"
$constN = "`$/unidad"
$constQ = 1
$constR = "Hortaliza"

# New rows (448..451): date 44855, two varieties x two qualities.
$newRows = @(
    @{ Row = 448; D = 44855; H = "Copenhague";    I = "Primera"; J = 400; K = 2500; L = 2500; M = 2500; O = "Región Metropolitana" },
    @{ Row = 449; D = 44855; H = "Copenhague";    I = "Segunda"; J = 300; K = 2300; L = 2300; M = 2300; O = "Región Metropolitana" },
    @{ Row = 450; D = 44855; H = "Crespo record"; I = "Primera"; J = 400; K = 2200; L = 2200; M = 2200; O = "Región Metropolitana" },
    @{ Row = 451; D = 44855; H = "Crespo record"; I = "Segunda"; J = 400; K = 2000; L = 2000; M = 2000; O = "Región Metropolitana" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $constN
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.M
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
}
